# Scheduled market-data refresh: update currentAveragePrice* / LevePrice*
# / LeveProfit* columns (H, I, J, K, L, M, N) across the eight Leve tables
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with freshly pulled values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H51").Value = 3333.1667   # currentAveragePrice row 51: was 4000
$ws.Range("J51").Value = 3333.1667   # currentAveragePriceHQ row 51: was 4000
$ws.Range("L51").Value = 3333.1667   # LevePriceHQ row 51: was 4000
$ws.Range("N51").Value = -4301.1667   # LeveProfitHQ row 51: was -4968

$ws.Range("H113").Value = 3259.6   # currentAveragePrice row 113: was 3369.7144
$ws.Range("I113").Value = 3366   # currentAveragePriceNQ row 113: was 3477.6
$ws.Range("K113").Value = 3366   # LevePriceNQ row 113: was 3477.6
$ws.Range("M113").Value = -112   # LeveProfitNQ row 113: was -223.5999999999999

$ws.Range("H132").Value = 4155.4375   # currentAveragePrice row 132: was 4332.467

$ws.Range("H137").Value = 2826.5557   # currentAveragePrice row 137: was 2963.4119
$ws.Range("I137").Value = 1249.875   # currentAveragePriceNQ row 137: was 1357
$ws.Range("K137").Value = 3749.625   # LevePriceNQ row 137: was 4071
$ws.Range("M137").Value = -1199.625   # LeveProfitNQ row 137: was -1521

$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 9096.134   # currentAveragePrice row 32: was 11338
$ws.Range("I32").Value = 9438.857   # currentAveragePriceNQ row 32: was 12042
$ws.Range("K32").Value = 9438.857   # LevePriceNQ row 32: was 12042
$ws.Range("M32").Value = -9151.857   # LeveProfitNQ row 32: was -11755

$ws.Range("H132").Value = 1568.5927   # currentAveragePrice row 132: was 1635.3334
$ws.Range("I132").Value = 1568.5927   # currentAveragePriceNQ row 132: was 1635.3334
$ws.Range("K132").Value = 4705.7781   # LevePriceNQ row 132: was 4906.0002
$ws.Range("M132").Value = -2175.7781   # LeveProfitNQ row 132: was -2376.0002

$ws = $wb.Worksheets("BSM")
$ws.Range("H86").Value = 3499.5   # currentAveragePrice row 86: was 3500
$ws.Range("I86").Value = 3499.5   # currentAveragePriceNQ row 86: was 3500
$ws.Range("K86").Value = 3499.5   # LevePriceNQ row 86: was 3500
$ws.Range("M86").Value = -2376.5   # LeveProfitNQ row 86: was -2377

$ws.Range("H89").Value = 3499.5   # currentAveragePrice row 89: was 3500
$ws.Range("I89").Value = 3499.5   # currentAveragePriceNQ row 89: was 3500
$ws.Range("K89").Value = 17497.5   # LevePriceNQ row 89: was 17500
$ws.Range("M89").Value = -11881.5   # LeveProfitNQ row 89: was -11884

$ws.Range("H107").Value = 1299.5   # currentAveragePrice row 107: was 1094.5
$ws.Range("I107").Value = 949.5   # currentAveragePriceNQ row 107: was 792.8333
$ws.Range("K107").Value = 949.5   # LevePriceNQ row 107: was 792.8333
$ws.Range("M107").Value = 970.5   # LeveProfitNQ row 107: was 1127.1667

$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 2503.4285   # currentAveragePrice row 31: was 2526.3416
$ws.Range("I31").Value = 2030.8966   # currentAveragePriceNQ row 31: was 2047.5714
$ws.Range("K31").Value = 2030.8966   # LevePriceNQ row 31: was 2047.5714
$ws.Range("M31").Value = -1735.8966   # LeveProfitNQ row 31: was -1752.5714

$ws.Range("H34").Value = 2503.4285   # currentAveragePrice row 34: was 2526.3416
$ws.Range("I34").Value = 2030.8966   # currentAveragePriceNQ row 34: was 2047.5714
$ws.Range("K34").Value = 2030.8966   # LevePriceNQ row 34: was 2047.5714
$ws.Range("M34").Value = -1828.8966   # LeveProfitNQ row 34: was -1845.5714

$ws.Range("H93").Value = 16500   # currentAveragePrice row 93: was 12332.333
$ws.Range("I93").Value = 16500   # currentAveragePriceNQ row 93: was 12332.333
$ws.Range("K93").Value = 16500   # LevePriceNQ row 93: was 12332.333
$ws.Range("M93").Value = -14628   # LeveProfitNQ row 93: was -10460.333

$ws.Range("H99").Value = 4896.8   # currentAveragePrice row 99: was 4966.7334
$ws.Range("I99").Value = 4813.909   # currentAveragePriceNQ row 99: was 5000.3
$ws.Range("J99").Value = 5124.75   # currentAveragePriceHQ row 99: was 4899.6
$ws.Range("K99").Value = 4813.909   # LevePriceNQ row 99: was 5000.3
$ws.Range("L99").Value = 5124.75   # LevePriceHQ row 99: was 4899.6
$ws.Range("M99").Value = -3315.909   # LeveProfitNQ row 99: was -3502.3
$ws.Range("N99").Value = -8120.75   # LeveProfitHQ row 99: was -7895.6

$ws.Range("H126").Value = 4896.8   # currentAveragePrice row 126: was 4966.7334
$ws.Range("I126").Value = 4813.909   # currentAveragePriceNQ row 126: was 5000.3
$ws.Range("J126").Value = 5124.75   # currentAveragePriceHQ row 126: was 4899.6
$ws.Range("K126").Value = 14441.727   # LevePriceNQ row 126: was 15000.9
$ws.Range("L126").Value = 15374.25   # LevePriceHQ row 126: was 14698.8
$ws.Range("M126").Value = -11971.727   # LeveProfitNQ row 126: was -12530.9
$ws.Range("N126").Value = -20314.25   # LeveProfitHQ row 126: was -19638.8

$ws.Range("H132").Value = 1028.5555   # currentAveragePrice row 132: was 1047.125
$ws.Range("I132").Value = 1028.5555   # currentAveragePriceNQ row 132: was 1047.125
$ws.Range("K132").Value = 3085.6665   # LevePriceNQ row 132: was 3141.375
$ws.Range("M132").Value = -555.6664999999998   # LeveProfitNQ row 132: was -611.375

$ws.Range("H134").Value = 2022.3846   # currentAveragePrice row 134: was 1527.4762
$ws.Range("I134").Value = 2024.25   # currentAveragePriceNQ row 134: was 1579.3334
$ws.Range("J134").Value = 2000   # currentAveragePriceHQ row 134: was 1216.3334
$ws.Range("K134").Value = 6072.75   # LevePriceNQ row 134: was 4738.0002
$ws.Range("L134").Value = 6000   # LevePriceHQ row 134: was 3649.0002
$ws.Range("M134").Value = -3537.75   # LeveProfitNQ row 134: was -2203.0002
$ws.Range("N134").Value = -11070   # LeveProfitHQ row 134: was -8719.0002

$ws = $wb.Worksheets("CUL")
$ws.Range("H11").Value = 678   # currentAveragePrice row 11: was 621.75
$ws.Range("I11").Value = 1006.6667   # currentAveragePriceNQ row 11: was 755.75
$ws.Range("K11").Value = 3020.0001   # LevePriceNQ row 11: was 2267.25
$ws.Range("M11").Value = -2880.0001   # LeveProfitNQ row 11: was -2127.25

$ws.Range("H113").Value = 1245.3334   # currentAveragePrice row 113: was 1197.5
$ws.Range("J113").Value = 1474.25   # currentAveragePriceHQ row 113: was 1363.7778
$ws.Range("L113").Value = 4422.75   # LevePriceHQ row 113: was 4091.3334
$ws.Range("N113").Value = -8762.75   # LeveProfitHQ row 113: was -8431.3334

$ws.Range("H117").Value = 749.8333   # currentAveragePrice row 117: was 753.5
$ws.Range("I117").Value = 754   # currentAveragePriceNQ row 117: was 755.3333
$ws.Range("J117").Value = 741.5   # currentAveragePriceHQ row 117: was 748
$ws.Range("K117").Value = 2262   # LevePriceNQ row 117: was 2265.9999
$ws.Range("L117").Value = 2224.5   # LevePriceHQ row 117: was 2244
$ws.Range("M117").Value = 1180   # LeveProfitNQ row 117: was 1176.0001
$ws.Range("N117").Value = -9108.5   # LeveProfitHQ row 117: was -9128

$ws.Range("H129").Value = 585   # currentAveragePrice row 129: was 622.2857
$ws.Range("J129").Value = 0   # currentAveragePriceHQ row 129: was 715.5
$ws.Range("L129").Value = 0   # LevePriceHQ row 129: was 2146.5
$ws.Range("N129").Value = ""   # LeveProfitHQ row 129: clear (was -12146.5)

$ws.Range("H131").Value = 1382   # currentAveragePrice row 131: was 1418
$ws.Range("I131").Value = 997.2857   # currentAveragePriceNQ row 131: was 999.1667
$ws.Range("K131").Value = 2991.8571   # LevePriceNQ row 131: was 2997.5001
$ws.Range("M131").Value = 2048.1429   # LeveProfitNQ row 131: was 2042.4999

$ws.Range("H139").Value = 2844.6667   # currentAveragePrice row 139: was 3131
$ws.Range("I139").Value = 2443.1428   # currentAveragePriceNQ row 139: was 2758
$ws.Range("K139").Value = 7329.428400000001   # LevePriceNQ row 139: was 8274
$ws.Range("M139").Value = -2189.428400000001   # LeveProfitNQ row 139: was -3134

$ws.Range("H140").Value = 530   # currentAveragePrice row 140: was 284
$ws.Range("I140").Value = 530   # currentAveragePriceNQ row 140: was 284
$ws.Range("K140").Value = 1590   # LevePriceNQ row 140: was 852
$ws.Range("M140").Value = 3590   # LeveProfitNQ row 140: was 4328

$ws = $wb.Worksheets("GSM")
$ws.Range("H7").Value = 34005000   # currentAveragePrice row 7: was 0
$ws.Range("I7").Value = 10002   # currentAveragePriceNQ row 7: was 0
$ws.Range("J7").Value = 68000000   # currentAveragePriceHQ row 7: was 0
$ws.Range("K7").Value = 10002   # LevePriceNQ row 7: was 0
$ws.Range("L7").Value = 68000000   # LevePriceHQ row 7: was 0
$ws.Range("M7").Value = -9890   # LeveProfitNQ row 7: new cell
$ws.Range("N7").Value = -68000224   # LeveProfitHQ row 7: new cell

$ws.Range("H8").Value = 34005000   # currentAveragePrice row 8: was 0
$ws.Range("I8").Value = 10002   # currentAveragePriceNQ row 8: was 0
$ws.Range("J8").Value = 68000000   # currentAveragePriceHQ row 8: was 0
$ws.Range("K8").Value = 10002   # LevePriceNQ row 8: was 0
$ws.Range("L8").Value = 68000000   # LevePriceHQ row 8: was 0
$ws.Range("M8").Value = -9863   # LeveProfitNQ row 8: new cell
$ws.Range("N8").Value = -68000278   # LeveProfitHQ row 8: new cell

$ws.Range("H39").Value = 0   # currentAveragePrice row 39: was 20000
$ws.Range("J39").Value = 0   # currentAveragePriceHQ row 39: was 20000
$ws.Range("L39").Value = ""   # LevePriceHQ row 39: clear (was 20000)
$ws.Range("N39").Value = 0   # LeveProfitHQ row 39: was -21064

$ws.Range("H97").Value = 1125.3334   # currentAveragePrice row 97: was 529.7857
$ws.Range("I97").Value = 719   # currentAveragePriceNQ row 97: was 256.54544
$ws.Range("K97").Value = 719   # LevePriceNQ row 97: was 256.54544
$ws.Range("M97").Value = -223   # LeveProfitNQ row 97: was 239.45456

$ws = $wb.Worksheets("LTW")
$ws.Range("H22").Value = 1088.4445   # currentAveragePrice row 22: was 1224
$ws.Range("J22").Value = 980   # currentAveragePriceHQ row 22: was 0
$ws.Range("L22").Value = 980   # LevePriceHQ row 22: was 0
$ws.Range("N22").Value = -1570   # LeveProfitHQ row 22: new cell

$ws.Range("H27").Value = 1088.4445   # currentAveragePrice row 27: was 1224
$ws.Range("J27").Value = 980   # currentAveragePriceHQ row 27: was 0
$ws.Range("L27").Value = 980   # LevePriceHQ row 27: was 0
$ws.Range("N27").Value = -1194   # LeveProfitHQ row 27: new cell

$ws.Range("H46").Value = 3694.95   # currentAveragePrice row 46: was 4230.6924
$ws.Range("I46").Value = 2612.5   # currentAveragePriceNQ row 46: was 2000
$ws.Range("K46").Value = 2612.5   # LevePriceNQ row 46: was 2000
$ws.Range("M46").Value = -2424.5   # LeveProfitNQ row 46: was -1812

$ws.Range("H61").Value = 451.5   # currentAveragePrice row 61: was 0
$ws.Range("I61").Value = 451.5   # currentAveragePriceNQ row 61: was 0
$ws.Range("K61").Value = 451.5   # LevePriceNQ row 61: was 0
$ws.Range("M61").Value = -249.5   # LeveProfitNQ row 61: new cell

$ws.Range("H113").Value = 451.5   # currentAveragePrice row 113: was 0
$ws.Range("I113").Value = 451.5   # currentAveragePriceNQ row 113: was 0
$ws.Range("K113").Value = 451.5   # LevePriceNQ row 113: was 0
$ws.Range("M113").Value = 1718.5   # LeveProfitNQ row 113: new cell

$ws.Range("H130").Value = 0   # currentAveragePrice row 130: was 66500.5
$ws.Range("J130").Value = 0   # currentAveragePriceHQ row 130: was 66500.5
$ws.Range("L130").Value = ""   # LevePriceHQ row 130: clear (was 66500.5)
$ws.Range("N130").Value = 0   # LeveProfitHQ row 130: was -76540.5

$ws.Range("H132").Value = 5021.316   # currentAveragePrice row 132: was 5189.5
$ws.Range("I132").Value = 5117   # currentAveragePriceNQ row 132: was 5741.6
$ws.Range("K132").Value = 15351   # LevePriceNQ row 132: was 17224.8
$ws.Range("M132").Value = -12821   # LeveProfitNQ row 132: was -14694.8

$ws.Range("H136").Value = 3500.5   # currentAveragePrice row 136: was 3504
$ws.Range("I136").Value = 3503.5   # currentAveragePriceNQ row 136: was 3504
$ws.Range("J136").Value = 3497.5   # currentAveragePriceHQ row 136: was 0
$ws.Range("K136").Value = 10510.5   # LevePriceNQ row 136: was 10512
$ws.Range("L136").Value = 10492.5   # LevePriceHQ row 136: was 0
$ws.Range("M136").Value = -7960.5   # LeveProfitNQ row 136: was -7962
$ws.Range("N136").Value = -15592.5   # LeveProfitHQ row 136: new cell

$ws = $wb.Worksheets("WVR")
$ws.Range("H62").Value = 3192.625   # currentAveragePrice row 62: was 4049
$ws.Range("I62").Value = 2923   # currentAveragePriceNQ row 62: was 4064.8333
$ws.Range("K62").Value = 2923   # LevePriceNQ row 62: was 4064.8333
$ws.Range("M62").Value = -2299   # LeveProfitNQ row 62: was -3440.8333

$ws.Range("H65").Value = 3192.625   # currentAveragePrice row 65: was 4049
$ws.Range("I65").Value = 2923   # currentAveragePriceNQ row 65: was 4064.8333
$ws.Range("K65").Value = 14615   # LevePriceNQ row 65: was 20324.1665
$ws.Range("M65").Value = -11495   # LeveProfitNQ row 65: was -17204.1665

$ws.Range("H75").Value = 90000   # currentAveragePrice row 75: was 5263
$ws.Range("I75").Value = 90000   # currentAveragePriceNQ row 75: was 5263
$ws.Range("K75").Value = 90000   # LevePriceNQ row 75: was 5263
$ws.Range("M75").Value = -89064   # LeveProfitNQ row 75: was -4327

$ws.Range("H78").Value = 90000   # currentAveragePrice row 78: was 5263
$ws.Range("I78").Value = 90000   # currentAveragePriceNQ row 78: was 5263
$ws.Range("K78").Value = 270000   # LevePriceNQ row 78: was 15789
$ws.Range("M78").Value = -265320   # LeveProfitNQ row 78: was -11109

$ws.Range("H81").Value = 885.8570999999999   # currentAveragePrice row 81: was 857.2857
$ws.Range("I81").Value = 885.8570999999999   # currentAveragePriceNQ row 81: was 857.2857
$ws.Range("K81").Value = 1771.7142   # LevePriceNQ row 81: was 1714.5714
$ws.Range("M81").Value = -710.7141999999999   # LeveProfitNQ row 81: was -653.5714

$ws.Range("H84").Value = 885.8570999999999   # currentAveragePrice row 84: was 857.2857
$ws.Range("I84").Value = 885.8570999999999   # currentAveragePriceNQ row 84: was 857.2857
$ws.Range("K84").Value = 8858.571   # LevePriceNQ row 84: was 8572.857
$ws.Range("M84").Value = -3554.571   # LeveProfitNQ row 84: was -3268.857

$ws.Range("H104").Value = 30370   # currentAveragePrice row 104: was 0
$ws.Range("J104").Value = 30370   # currentAveragePriceHQ row 104: was 0
$ws.Range("L104").Value = 30370   # LevePriceHQ row 104: was 0
$ws.Range("N104").Value = -37358   # LeveProfitHQ row 104: new cell

$ws.Range("H108").Value = 55000   # currentAveragePrice row 108: was 0
$ws.Range("J108").Value = 55000   # currentAveragePriceHQ row 108: was 0
$ws.Range("L108").Value = 55000   # LevePriceHQ row 108: was 0
$ws.Range("N108").Value = -62680   # LeveProfitHQ row 108: new cell

$ws.Range("H109").Value = 100376   # currentAveragePrice row 109: was 100376.75
$ws.Range("J109").Value = 100376   # currentAveragePriceHQ row 109: was 100376.75
$ws.Range("L109").Value = 100376   # LevePriceHQ row 109: was 100376.75
$ws.Range("N109").Value = -103150   # LeveProfitHQ row 109: was -103150.75

$ws.Range("H126").Value = 1248.4286   # currentAveragePrice row 126: was 1497.3334
$ws.Range("I126").Value = 1248.4286   # currentAveragePriceNQ row 126: was 1497.3334
$ws.Range("K126").Value = 3745.2858   # LevePriceNQ row 126: was 4492.0002
$ws.Range("M126").Value = -1275.2858   # LeveProfitNQ row 126: was -2022.0002

$ws.Range("H132").Value = 1460.9131   # currentAveragePrice row 132: was 1560.2727
$ws.Range("I132").Value = 1305   # currentAveragePriceNQ row 132: was 1406.5
$ws.Range("K132").Value = 3915   # LevePriceNQ row 132: was 4219.5
$ws.Range("M132").Value = -1385   # LeveProfitNQ row 132: was -1689.5

$ws.Range("H136").Value = 11599   # currentAveragePrice row 136: was 11998.5
$ws.Range("I136").Value = 11758.8   # currentAveragePriceNQ row 136: was 11998.5
$ws.Range("J136").Value = 10800   # currentAveragePriceHQ row 136: was 0
$ws.Range("K136").Value = 35276.39999999999   # LevePriceNQ row 136: was 35995.5
$ws.Range("L136").Value = 32400   # LevePriceHQ row 136: was 0
$ws.Range("M136").Value = -32726.39999999999   # LeveProfitNQ row 136: was -33445.5
$ws.Range("N136").Value = -37500   # LeveProfitHQ row 136: new cell
